$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 previously held the "122434 / 07:41:18" log entry, row 3 held the
# "121234 / 07:41:20" entry. The edit removes row 2's original entry and
# keeps/shifts the row 3 entry up into row 2, so update A2/D2 in place
# (the rest of row 2 - Subject/Log Date/Type/User - already match row 3)
# and then delete the now-duplicate row 3.

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "121234"
$ws.Range("D2").Value = "07:41:20"

$ws.Range("A3:F3").EntireRow.Delete()
